$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion summary text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$conversionText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.11 = 7726.32 pesos`n✅ 7726.32 pesos = 2.1 = 934.69 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $conversionText

# --- Sheet "tasas": update the rate cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("O10").Value = 3670
$wsTasas.Range("N12").Value = 3685
$wsTasas.Range("O12").Value = 445.79
